$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Revert phone number prefixes in column A (tel) for rows 5-25.
$ws.Range("A5").Value = "333333333333"
$ws.Range("A6").Value = "444444444444"
$ws.Range("A7").Value = "555555555555"
$ws.Range("A11").Value = "222223333333"
$ws.Range("A12").Value = "444447777777"
$ws.Range("A13").Value = "555555987654"
$ws.Range("A14").Value = "123438459832"
$ws.Range("A17").Value = "333333320098"
$ws.Range("A18").Value = "488829844444"
$ws.Range("A19").Value = "555998055555"
$ws.Range("A23").Value = "222229933333"
$ws.Range("A24").Value = "447777777777"
$ws.Range("A25").Value = "555553437654"

# Update the active cell selection on the sheet.
$ws.Range("I10").Select()
